$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-30 Tuesday" "2024-07-31 Wednesday"

Replace-Text "40÷6=" "28÷9="
Replace-Text "64÷3=" "39÷5="
Replace-Text "52÷2=" "48÷3="
Replace-Text "67÷3=" "98÷5="
Replace-Text "84÷2=" "47÷5="
Replace-Text "47÷2=" "98÷8="
Replace-Text "75÷8=" "13÷6="
Replace-Text "90÷7=" "16÷7="
Replace-Text "86÷6=" "18÷4="
Replace-Text "25÷6=" "43÷8="
Replace-Text "12÷5=" "77÷6="
Replace-Text "74÷7=" "25÷5="
Replace-Text "64÷4=" "93÷9="
Replace-Text "84÷4=" "92÷6="
Replace-Text "43÷2=" "21÷4="
Replace-Text "31÷9=" "86÷4="
Replace-Text "13÷3=" "35÷2="
Replace-Text "91÷5=" "58÷2="
Replace-Text "65÷2=" "53÷8="
Replace-Text "12÷3=" "47÷5="
Replace-Text "20÷2=" "46÷6="
Replace-Text "69÷2=" "56÷7="
Replace-Text "68÷5=" "28÷7="
Replace-Text "33÷7=" "74÷6="
Replace-Text "69÷3=" "30÷8="

Write-Output "Done"
